$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 71.541692
$ws.Range("N2").Value = 214.625076
$ws.Range("O2").Value = 0.6133071420247926
$ws.Range("P2").Value = 0.6133071420247926
$ws.Range("Q2").Value = 4.025340994841333
$ws.Range("R2").Value = 36.228068953572
$ws.Range("S2").Value = 0.1673486831912561
$ws.Range("T2").Value = 0.1673486831912561
$ws.Range("O3").Value = 0.08457024278578675
$ws.Range("P3").Value = 0.08457024278578675
$ws.Range("S3").Value = 0.02307607036931281
$ws.Range("T3").Value = 0.02307607036931281
$ws.Range("M4").Value = 35.05835333333334
$ws.Range("N4").Value = 105.17506
$ws.Range("O4").Value = 0.3005455684073286
$ws.Range("P4").Value = 0.3005455684073286
$ws.Range("Q4").Value = 1.972581622535556
$ws.Range("R4").Value = 17.75323460282
$ws.Range("S4").Value = 0.08200769510996636
$ws.Range("T4").Value = 0.08200769510996636
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.183961
$ws.Range("N5").Value = 0.551883
$ws.Range("O5").Value = 0.001577046782092083
$ws.Range("P5").Value = 0.001577046782092083
$ws.Range("Q5").Value = 0.01035068830566667
$ws.Range("R5").Value = 0.09315619475100001
$ws.Range("S5").Value = 0.000430317347100858
$ws.Range("T5").Value = 0.000430317347100858
$ws.Range("M6").Value = 71.541692
$ws.Range("N6").Value = 214.625076
$ws.Range("O6").Value = 0.6133071420247926
$ws.Range("P6").Value = 0.6133071420247926
$ws.Range("Q6").Value = 10.72691360401867
$ws.Range("R6").Value = 96.542222436168
$ws.Range("S6").Value = 0.4459584588335365
$ws.Range("T6").Value = 0.4459584588335365
$ws.Range("O7").Value = 0.08457024278578675
$ws.Range("P7").Value = 0.08457024278578675
$ws.Range("S7").Value = 0.06149417241647393
$ws.Range("T7").Value = 0.06149417241647393
$ws.Range("M8").Value = 35.05835333333334
$ws.Range("N8").Value = 105.17506
$ws.Range("O8").Value = 0.3005455684073286
$ws.Range("P8").Value = 0.3005455684073286
$ws.Range("Q8").Value = 5.256626126564445
$ws.Range("R8").Value = 47.30963513908001
$ws.Range("S8").Value = 0.2185378732973622
$ws.Range("T8").Value = 0.2185378732973622
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.183961
$ws.Range("N9").Value = 0.551883
$ws.Range("O9").Value = 0.001577046782092083
$ws.Range("P9").Value = 0.001577046782092083
$ws.Range("Q9").Value = 0.02758298969933334
$ws.Range("R9").Value = 0.248246907294
$ws.Range("S9").Value = 0.001146729434991225
$ws.Range("T9").Value = 0.001146729434991225